$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $r = $d.Content
    $ok = $r.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Find.Execute failed to find: $findText"
    }
    # Assign .Text directly (rather than using Find.Execute's replace argument)
    # so Word's smart-quote AutoCorrect does not mangle straight apostrophes.
    $r.Text = $replaceText
}

function Append-After($findText, $appendText) {
    $r = $d.Content
    $ok = $r.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Find.Execute failed to find (for append): $findText"
    }
    $r.Collapse(0)
    $r.InsertAfter($appendText)
}

# --- Title ---
Replace-Text "Ascend to Knowledge, Redefining Horizon" "Exploring the Wonders of Biology: The Science of Life"

# --- Author ---
Replace-Text "William Sides" "Alex Smith"

# --- Email line (runs collapse naturally) ---
Replace-Text "research@gnail.org" "edu"
Replace-Text "academic" "alexsmith@bioworld"

# --- Body paragraph 1 ---
Replace-Text "In the boundless labyrinth of knowledge, the ceaseless quest for truth and understanding propels us forward, transcending limits and shifting paradigms" "Biology, the study of life, unravels the complexities and marvels of living organisms"

Replace-Text " We set sail on the sea of inquiry, our minds eager to encounter new horizons, to unravel enigmas, and to reveal the intricate tapestry of the universe" " It delves into the intricate mechanisms that govern our bodies, the beauty of ecosystems, and the interactions between organisms and their environment"

Replace-Text " Through the gateway of science, technology, and scholarship, we uncover marvels of the natural world, unravel the complexities of the cosmos, and glimpse the infinite facets of human endeavor" " As we explore the wonders of biology, we embark on a journey of discovery, uncovering the secrets of our existence and unlocking the mysteries of nature"
Append-After " As we explore the wonders of biology, we embark on a journey of discovery, uncovering the secrets of our existence and unlocking the mysteries of nature" ". From the microscopic world of cells to the vastness of ecosystems, biology captivates us with its elegance and interconnectedness"

Replace-Text "As we delve into the depths of knowledge, we encounter mysteries and phenomena that challenge our convictions and reshape our perceptions" "As we delve deeper into biology, we uncover the fundamental principles that govern life"

Replace-Text " Innovations emerge from the laboratories, pushing the frontiers of human capability, and unlocking doors to realms once deemed unreachable" " We unravel the secrets of DNA, the molecule that holds the blueprints for all living organisms"

Replace-Text " From the subtle interactions of subatomic particles to the vast expanses of the cosmos, we embark on a voyage of discovery, seeking to comprehend the fundamental forces that shape our existence" " We explore the intricacies of cells, the building blocks of life, and the remarkable adaptations that allow organisms to thrive in diverse environments"
Append-After " We explore the intricacies of cells, the building blocks of life, and the remarkable adaptations that allow organisms to thrive in diverse environments" ". Biology reveals the interconnectedness of life, demonstrating how organisms depend on each other in intricate ecosystems, forming a harmonious web of interactions"

Replace-Text "Each discipline, like a mosaic tile, adds a fragment of insight to the grand pattern of comprehension" "Biology's impact extends far beyond the laboratory"

Replace-Text " The tapestry of knowledge is intricately woven, with threads of logic, creativity, and empirical investigation intertwining" " It has revolutionized medicine, leading to life-saving drugs, therapies, and treatments"

Replace-Text " We strive to grasp the essence of reality, to forge connections between seemingly disparate fields, and to glean profound truths from the careful examination of evidence" " It has transformed agriculture, increasing crop yields and improving food security"
Append-After " It has transformed agriculture, increasing crop yields and improving food security" ". Biology empowers us to understand and address environmental challenges, promoting sustainability and conservation"
Append-After "Biology empowers us to understand and address environmental challenges, promoting sustainability and conservation" ". As we unravel the mysteries of life, we gain a profound appreciation for the beauty and complexity of our planet"

# --- Summary paragraph ---
Replace-Text "Our insatiable need to understand the world around us ignites the journey of knowledge acquisition" "Biology, the study of life, unveils the intricate workings of living organisms, the elegance of ecosystems, and the fundamental principles that govern life"

Replace-Text " In this odyssey of learning, we embark on paths that lead to unexpected destinations, challenging our understanding and expanding our horizons" " It uncovers the secrets of DNA and cells, revealing the interconnectedness of life and the remarkable adaptations that allow organisms to thrive"

Replace-Text " We embrace uncertainty as a catalyst for growth, recognizing that knowledge is not a static entity but a dynamic process of exploration, discovery, and connection" " Biology has revolutionized medicine, agriculture, and our understanding of the environment, empowering us to address global challenges and foster sustainability"
Append-After " Biology has revolutionized medicine, agriculture, and our understanding of the environment, empowering us to address global challenges and foster sustainability" ". As we delve into the wonders of biology, we unlock the mysteries of life and gain a profound appreciation for the complexity and beauty of our planet"

# --- Add trailing empty paragraph before the end of the document body ---
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
